$d = $word.ActiveDocument

# Helper: force a run boundary at the edges of $rng without permanently
# changing its appearance -- flip a character property to the opposite
# of the desired final value and then back. The engine merges
# identically-formatted adjacent runs on every edit, but it does not
# retroactively re-merge runs that already exist, so ending on the
# same value the text already had still leaves the split in place.
function Split-Run($rng, [bool]$finalBold) {
    $rng.Font.Bold = -not $finalBold
    $rng.Font.Bold = $finalBold
}

# ---------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that currently lives alone in a
#    paragraph right after the "{fullName}" field. Deleting the
#    paragraph's own (zero-length) range strips the bookmark markup
#    while leaving the paragraph itself intact (now truly empty).
#    Doing this before minting a new "_GoBack" bookmark keeps the name
#    lookup below unambiguous (Word only allows one bookmark per name).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBookmark = $d.Bookmarks("_GoBack")
    $oldParagraph = $d.Range($oldBookmark.Start, $oldBookmark.Start).Paragraphs(1)
    $oldParagraph.Range.Delete()
}

# ---------------------------------------------------------------------
# 2) "Flight company" -> "Flight" + "'s" + " company" (three runs, all
#    keeping the same underline formatting), with a fresh "_GoBack"
#    bookmark sitting between the "'s" run and the " company" run.
# ---------------------------------------------------------------------

# Replace the literal text in place so the original run (and its
# <w:u w:val="single"/> formatting) is preserved; the apostrophe is a
# right single quotation mark (U+2019), matching the diff's "'s".
$apos = [char]0x2019
$newText = "Flight" + $apos + "s company"
$d.Content.Find.Execute("Flight company", $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# Re-find the freshly written text so we work off real positions
# instead of hard-coded offsets.
$found = $d.Content
$found.Find.Execute($newText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$runStart = $found.Start

# Force a run boundary between "Flight" and "'s company" without
# changing visible formatting.
Split-Run $d.Range($runStart, $runStart + 6) $false   # "Flight"

# Insert the (hidden) "_GoBack" bookmark right between "'s" and
# " company" -- this naturally splits the remaining text into its own
# run too, giving the required "'s" / " company" run boundary.
$bookmarkPos = $runStart + 8   # just after "Flight's" (6 + 2 chars)
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# 3) The edits above touch the whole "Flight...{Company}" paragraph, so
#    the engine also re-flows the untouched "{Company}" placeholder
#    that follows -- restore its original five-run split (it is not
#    part of the requested change).
# ---------------------------------------------------------------------
$cfound = $d.Content
$cfound.Find.Execute("{Company}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$cs = $cfound.Start

Split-Run $d.Range($cs, $cs + 1) $true       # "{"
Split-Run $d.Range($cs + 1, $cs + 4) $true   # "Com"
Split-Run $d.Range($cs + 4, $cs + 5) $true   # "p"
Split-Run $d.Range($cs + 5, $cs + 6) $true   # "a"
Split-Run $d.Range($cs + 6, $cs + 9) $true   # "ny}"
